# New staging templates generated:
# Add "EndDate" as a new first data column and "StartDate" as a new last
# data column to the header row, shifting the existing header labels
# (OrganizationBusinessKey, OrganizationPersonRole_ID, PersonBusinessKey,
# RoleBusinessKey) one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the bold+underline formatting already used by the header row
# (row 2) so the two new header cells match it exactly.
$headerFont = $ws.Range("B2").Font

# Shift the existing header values one column to the right (B2:E2),
# freeing up A2 for the new "EndDate" column and making room for the
# new "StartDate" column in F2.
$ws.Range("E2").Value = $ws.Range("D2").Value2
$ws.Range("D2").Value = $ws.Range("C2").Value2
$ws.Range("C2").Value = $ws.Range("B2").Value2
$ws.Range("B2").Value = $ws.Range("A2").Value2

# Set the new header labels.
$ws.Range("A2").Value = "EndDate"
$ws.Range("F2").Value = "StartDate"

# Make sure every header cell (including the two new ones) keeps the
# bold + underline formatting used throughout row 2.
$headerRow = $ws.Range("A2:F2")
$headerRow.Font.Bold = $headerFont.Bold
$headerRow.Font.Underline = $headerFont.Underline
